# A new weekly price record was added to the top of this "Espinaca" /
# "Vega Modelo de Temuco" sub-table. In the worksheet this shows up as a
# brand-new row inserted right above the existing row 261 (the table is
# sorted most-recent-first and row 261 was, until now, the most recent
# "$/docena de atados" entry). Inserting the row pushes every following
# record down by one (old 261 -> 262, ... old 294 -> 295), and the new
# row at 261 is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldRow = 261
$lastCol = 18   # columns A..R

# Insert a blank row above row 261; everything below shifts down one row.
$ws.Rows($oldRow).Insert()

# After the insert, the record that used to live at row 261 now lives at
# row 262. Re-populate the freshly inserted (now blank) row 261 with the
# same reference data (market, region, product, quality, unit, origin,
# etc.) before overwriting the two cells that actually carry this week's
# new numbers (date + volume).
$sourceRow = $oldRow + 1
for ($col = 1; $col -le $lastCol; $col++) {
    $value = $ws.Cells.Item($sourceRow, $col).Value()
    $ws.Cells.Item($oldRow, $col).Value = $value
}

# New weekly observation: Fecha (D) and Volumen (J).
$ws.Cells.Item($oldRow, 4).Value = 45131
$ws.Cells.Item($oldRow, 10).Value = 80
